$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so decimal-looking
# values like "0.9934" are not coerced into floating point numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.206.46'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '1.875.72'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('D4').Value = '0.9934'
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').Value = '236.49'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').Value = '0.9939'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('D7').Value = '0.4656'
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('D8').Value = '0.2828'
$ws.Range('E8').Value = '  +3.29%  '
$ws.Range('D9').Value = '0.06513'
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').Value = '108.20'
$ws.Range('E10').Value = '  +28.10%  '
$ws.Range('D11').Value = '18.71'
$ws.Range('E11').Value = '  +6.21%  '
$ws.Range('D12').Value = '1.849.10'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '0.07513'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '5.037'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '313.39'
$ws.Range('E15').Value = '  +29.66%  '
$ws.Range('D16').Value = '0.6330'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '30.195.63'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '0.9932'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').Value = '0.000007454'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '2.084.87'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').Value = '0.9935'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = '5.056'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').Value = '6.220'
$ws.Range('E24').Value = '  +4.39%  '
$ws.Range('D25').Value = '9.179'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '164.79'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('D27').Value = '20.13'
$ws.Range('E27').Value = '  +11.95%  '
$ws.Range('D28').Value = '1.984'
$ws.Range('E28').Value = '  +5.65%  '
$ws.Range('D29').Value = '0.1080'
$ws.Range('E29').Value = '  +5.98%  '
$ws.Range('D30').Value = '1.331'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('D31').Value = '4.048'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '3.882'
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('D33').Value = '0.04918'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').Value = '0.7433'
$ws.Range('E34').Value = '  +5.49%  '
$ws.Range('D35').Value = '1.135'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').Value = '2.701'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').Value = '0.01928'
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('D38').Value = '2.660'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').Value = '1.991'
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('D40').Value = '0.8649'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').Value = '106.97'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').Value = '5.738'
$ws.Range('E42').Value = '  +4.42%  '
$ws.Range('D43').Value = '0.9936'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').Value = '0.4102'
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').Value = '66.81'
$ws.Range('E45').Value = '  +7.39%  '
$ws.Range('D46').Value = '7.150'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '9.215'
$ws.Range('E47').Value = '  +8.60%  '
$ws.Range('D48').Value = '0.1201'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = '34.16'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('D50').Value = '0.05553'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').Value = '0.3759'
$ws.Range('E51').Value = '  +2.54%  '
